# Correct the mislabeled "Plain_English" header (underscore) to the proper
# "Plain English" (space) wording across every lookup/code table sheet that
# uses it: VehFuel_to_Code, Tech_to_Code, Dem_to_Code and Fuel_to_Code.

$wb = $excel.ActiveWorkbook

# VehFuel_to_Code : header in column C
$wsVehFuel = $wb.Worksheets.Item("VehFuel_to_Code")
$wsVehFuel.Range("C1").Value = "Plain English"
$wsVehFuel.Activate()
$wsVehFuel.Range("C1").Select()

# Tech_to_Code : header in column C
$wsTech = $wb.Worksheets.Item("Tech_to_Code")
$wsTech.Range("C1").Value = "Plain English"
$wsTech.Activate()
$wsTech.Range("C1").Select()

# Dem_to_Code : header in column B
$wsDem = $wb.Worksheets.Item("Dem_to_Code")
$wsDem.Range("B1").Value = "Plain English"
$wsDem.Activate()
$wsDem.Range("B1").Select()

# Fuel_to_Code : header in column C
$wsFuel = $wb.Worksheets.Item("Fuel_to_Code")
$wsFuel.Range("C1").Value = "Plain English"

# Leave the workbook focused on the sheet where the correction was made
# last, with the edited cell selected, matching the author's final
# editing state.
$wsFuel.Activate()
$wsFuel.Range("C1").Select()
